# Add a new "14-jul" column (AA) to Sheet1, mirroring the existing "13-jul"
# (Z) column with a handful of values tweaked, and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new date column (text, same style as the other date headers).
$ws.Range("AA1").Value = "14-jul"
$ws.Range("AA1").NumberFormat = "@"

# New column values (row 2..11), based on the previous "13-jul" (Z) column
# with a handful of values bumped/trimmed. Match the existing numeric
# formatting (centered integers) used by the rest of the data rows.
$ws.Range("AA2:AA11").HorizontalAlignment = -4108
$ws.Range("AA2:AA11").NumberFormat = "0"

$ws.Range("AA2").Value = 13
$ws.Range("AA3").Value = 20
$ws.Range("AA4").Value = 7
$ws.Range("AA5").Value = 9
$ws.Range("AA6").Value = 12
$ws.Range("AA7").Value = 16
$ws.Range("AA8").Value = 19
$ws.Range("AA9").Value = 13
$ws.Range("AA10").Value = 18
$ws.Range("AA11").Value = 22

# Move the active selection as recorded after the edit.
$ws.Range("AB8").Select()
